$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  @(-0.23336480490412725, 0.23318161450306718),
  @(-0.1566683690228956, 0.15628220052232589),
  @(-0.10657863379354637, 0.10630016582479485),
  @(-0.098300165849309451, 0.097823991695531021),
  @(-0.094823991709375832, 0.093204215894395759),
  @(-0.052730420987572302, 0.052074242979980312),
  @(-0.042079296499086816, 0.041929752581006063),
  @(-0.031929752616959739, 0.031648831022116131),
  @(-0.029648831040187229, 0.029412813384205094),
  @(-0.027412813403669745, 0.027396317957608218),
  @(-0.024396317979674897, 0.024369395781139502),
  @(-0.020869395804783419, 0.02067324902097134),
  @(-0.017173249046033234, 0.017083702057511196),
  @(-0.0090837020938918656, 0.0090542450537069996),
  @(-0.0080542450734135684, 0.008035204179507538),
  @(-0.0060352042018951835, 0.0060037362487341817),
  @(-0.0040037362714890889, 0.0039999999722892809),
  @(-0.053608152945422205, 0.053521073993412216),
  @(-0.04952107400406236, 0.048903932313542064),
  @(-0.044903932327367002, 0.044733375541394338),
  @(-0.040733375555944917, 0.040479060767984798),
  @(-0.045703566742153257, 0.045492572594470815),
  @(-0.040492572611036337, 0.040097820136409723),
  @(-0.02009782019066364, 0.019999999945055968),
  @(-0.038473794921726423, 0.038453613385689422),
  @(-0.035953613402066154, 0.035930727830820075),
  @(-0.033430727847514774, 0.033313059434276759),
  @(-0.031313059451212766, 0.031244862122660422),
  @(-0.02424486215271493, 0.024233273520844811),
  @(-0.021165168522453204, 0.021023377578856639),
  @(-0.014023377610229204, 0.014001542812078327),
  @(-0.0040015428508564099, 0.0039999999756759053)
)

for ($i = 0; $i -lt $values.Count; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $values[$i][0]
  $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

$ws.Columns.Item(2).ColumnWidth = 13.8333333
